# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-10-06 (serial 45205) to 2023-10-07 (serial 45206).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$newDate = Get-Date -Year 2023 -Month 10 -Day 7 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

for ($row = 2; $row -le 163; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value -ne $null) {
        $cell.Value = $newDate
    }
}
